$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1) updates to column F ("想去人数")
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 139
$ws1.Range("F6").Value = 1023
$ws1.Range("F7").Value = 2189
$ws1.Range("F9").Value = 1098
$ws1.Range("F10").Value = 601
$ws1.Range("F18").Value = 1572
$ws1.Range("F19").Value = 623
$ws1.Range("F21").Value = 596
$ws1.Range("F22").Value = 12198
$ws1.Range("F23").Value = 12229
$ws1.Range("F30").Value = 1917

# Sheet "全部类型" (sheet4) updates to column F ("想去人数")
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 139
$ws4.Range("F7").Value = 1023
$ws4.Range("F8").Value = 2189
$ws4.Range("F10").Value = 1098
$ws4.Range("F11").Value = 601
$ws4.Range("F22").Value = 1572
$ws4.Range("F23").Value = 623
$ws4.Range("F25").Value = 596
$ws4.Range("F26").Value = 12198
$ws4.Range("F27").Value = 12229
$ws4.Range("F34").Value = 1917
